# Update countries & provincias Spain
#
# Source data got re-scraped: a handful of countries (Colombia, Ruanda,
# Jamaica) received fresh case counts and moved up one row (ahead of the
# neighbour they used to trail), which bumps that neighbour's old numbers
# down a row; a few other countries got their counters refreshed in place;
# and the "last updated" footer timestamp advanced from 00:22 to 00:52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($row, $name, $vals)
    if ($name) { $ws.Cells.Item($row, 1).Value = $name }
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
    $ws.Cells.Item($row, 8).Value = $vals[6]
}

# Footer timestamp (row 1, column A)
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 00:52"

# Estados Unidos (row 4) - refreshed counters
Set-Row 4 $null @(530384, 27508, 29444, 480427, 11315, 1766, 20513)

# Colombia moves ahead of Tailandia/Ucrania with new counters; the two
# displaced countries keep their previous figures, shifted down one row.
Set-Row 49 "Colombia"  @(2709, 236, 214, 2395, 92, 20, 100)
Set-Row 50 "Tailandia" @(2518, 45, 1135, 1348, 61, 2, 35)
Set-Row 51 "Ucrania"   @(2511, 308, 79, 2359, 45, 4, 73)

# Ruanda moves ahead of Camboya with new counters; Camboya keeps its
# previous figures, shifted down one row.
Set-Row 128 "Ruanda"  @(120, 2, 18, 102, 0, 0, 0)
Set-Row 129 "Camboya" @(120, 1, 75, 45, 1, 0, 0)

# Jamaica moves ahead of Barbados with new counters; Barbados keeps its
# previous figures, shifted down one row.
Set-Row 140 "Jamaica"  @(69, 6, 13, 52, 0, 0, 4)
Set-Row 141 "Barbados" @(68, 1, 11, 53, 4, 0, 4)

# Guyana (row 152) - refreshed counters (D, F, G unchanged)
$ws.Cells.Item(152, 2).Value = 45
$ws.Cells.Item(152, 3).Value = 8
$ws.Cells.Item(152, 5).Value = 31
